$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revert the header row back to the original Ukrainian labels (was
# temporarily switched to English identifier-style headers).
$ws.Range("A1").Value = "Рік"
$ws.Range("B1").Value = "Об'єкт"
$ws.Range("C1").Value = "Вид діяльності"
$ws.Range("D1").Value = "Місцезнаходження"
$ws.Range("E1").Value = "Оксид азоту, т/рік**"
$ws.Range("F1").Value = "Cірки діоксид, т/рік**"
$ws.Range("G1").Value = "Оксид вуглецю, т/рік**"
$ws.Range("H1").Value = "Мікрочастинки та волокна, т/рік**"
$ws.Range("I1").Value = "Всього, т/рік**"

# Restore the active selection to K5 (and drop the prior frozen
# top-left scroll position on C1).
$ws.Range("K5").Select()
